$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1727.7778
$ws.Range("I40").Value = 800
$ws.Range("J40").Value = 1843.75
$ws.Range("K40").Value = 800
$ws.Range("L40").Value = 1843.75
$ws.Range("M40").Value = -625
$ws.Range("N40").Value = -2193.75
# Row 75
$ws.Range("H75").Value = 14999.5
$ws.Range("J75").Value = 14999.5
$ws.Range("L75").Value = 14999.5
$ws.Range("N75").Value = -16871.5
# Row 78
$ws.Range("H78").Value = 14999.5
$ws.Range("J78").Value = 14999.5
$ws.Range("L78").Value = 44998.5
$ws.Range("N78").Value = -54358.5
# Row 86
$ws.Range("H86").Value = 6246.615
$ws.Range("I86").Value = 5240.2
$ws.Range("J86").Value = 9601.333000000001
$ws.Range("K86").Value = 5240.2
$ws.Range("L86").Value = 9601.333000000001
$ws.Range("M86").Value = -4117.2
$ws.Range("N86").Value = -11847.333
# Row 89
$ws.Range("H89").Value = 6246.615
$ws.Range("I89").Value = 5240.2
$ws.Range("J89").Value = 9601.333000000001
$ws.Range("K89").Value = 26201
$ws.Range("L89").Value = 48006.665
$ws.Range("M89").Value = -20585
$ws.Range("N89").Value = -59238.665
# Row 134
$ws.Range("H134").Value = 60780
$ws.Range("J134").Value = 60780
$ws.Range("L134").Value = 60780
$ws.Range("N134").Value = -70920
# Row 137
$ws.Range("H137").Value = 1264
$ws.Range("I137").Value = 1091.1818
$ws.Range("J137").Value = 1609.6364
$ws.Range("K137").Value = 3273.5454
$ws.Range("L137").Value = 4828.9092
$ws.Range("M137").Value = -723.5454
$ws.Range("N137").Value = -9928.9092

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3537.8728
$ws.Range("I32").Value = 3597.0889
$ws.Range("J32").Value = 3271.4
$ws.Range("K32").Value = 3597.0889
$ws.Range("L32").Value = 3271.4
$ws.Range("M32").Value = -3310.0889
$ws.Range("N32").Value = -3845.4
# Row 74
$ws.Range("H74").Value = 1114.4333
$ws.Range("I74").Value = 746.7727
$ws.Range("K74").Value = 746.7727
$ws.Range("M74").Value = 127.2273
# Row 77
$ws.Range("H77").Value = 1114.4333
$ws.Range("I77").Value = 746.7727
$ws.Range("K77").Value = 3733.8635
$ws.Range("M77").Value = 634.1365000000001
# Row 88
$ws.Range("H88").Value = 2659.8635
$ws.Range("I88").Value = 2068.3333
$ws.Range("J88").Value = 2881.6875
$ws.Range("K88").Value = 2068.3333
$ws.Range("L88").Value = 2881.6875
$ws.Range("M88").Value = -1662.3333
$ws.Range("N88").Value = -3693.6875
# Row 91
$ws.Range("H91").Value = 2659.8635
$ws.Range("I91").Value = 2068.3333
$ws.Range("J91").Value = 2881.6875
$ws.Range("K91").Value = 2068.3333
$ws.Range("L91").Value = 2881.6875
$ws.Range("M91").Value = -664.3332999999998
$ws.Range("N91").Value = -5689.6875

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3432.25
$ws.Range("I86").Value = 4367.5264
$ws.Range("J86").Value = 2065.3076
$ws.Range("K86").Value = 4367.5264
$ws.Range("L86").Value = 2065.3076
$ws.Range("M86").Value = -3244.5264
$ws.Range("N86").Value = -4311.3076
# Row 89
$ws.Range("H89").Value = 3432.25
$ws.Range("I89").Value = 4367.5264
$ws.Range("J89").Value = 2065.3076
$ws.Range("K89").Value = 21837.632
$ws.Range("L89").Value = 10326.538
$ws.Range("M89").Value = -16221.632
$ws.Range("N89").Value = -21558.538

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1416
$ws.Range("I31").Value = 1272
$ws.Range("K31").Value = 1272
$ws.Range("M31").Value = -977
# Row 34
$ws.Range("H34").Value = 1416
$ws.Range("I34").Value = 1272
$ws.Range("K34").Value = 1272
$ws.Range("M34").Value = -1070
# Row 75
$ws.Range("H75").Value = 20260
$ws.Range("J75").Value = 20260
$ws.Range("L75").Value = 20260
$ws.Range("N75").Value = -22256
# Row 78
$ws.Range("H78").Value = 20260
$ws.Range("J78").Value = 20260
$ws.Range("L78").Value = 60780
$ws.Range("N78").Value = -70764
# Row 132
$ws.Range("H132").Value = 4881.3823
$ws.Range("I132").Value = 5550.304
$ws.Range("J132").Value = 3482.7273
$ws.Range("K132").Value = 16650.912
$ws.Range("L132").Value = 10448.1819
$ws.Range("M132").Value = -14120.912
$ws.Range("N132").Value = -15508.1819

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 1045875.2
$ws.Range("I4").Value = 583029.7
$ws.Range("J4").Value = 1200157
$ws.Range("K4").Value = 1749089.1
$ws.Range("L4").Value = 3600471
$ws.Range("M4").Value = -1748977.1
$ws.Range("N4").Value = -3600695
# Row 131
$ws.Range("H131").Value = 14926593
$ws.Range("J131").Value = 1317.2034
$ws.Range("L131").Value = 3951.6102
$ws.Range("N131").Value = -14031.6102
# Row 133
$ws.Range("H133").Value = 3915.5334
$ws.Range("J133").Value = 4367.154
$ws.Range("L133").Value = 13101.462
$ws.Range("N133").Value = -23221.462

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 24
$ws.Range("H24").Value = 1115322.9
$ws.Range("I24").Value = 5001953
$ws.Range("J24").Value = 4857.143
$ws.Range("K24").Value = 5001953
$ws.Range("L24").Value = 4857.143
$ws.Range("M24").Value = -5001780
$ws.Range("N24").Value = -5203.143
# Row 34
$ws.Range("H34").Value = 29999
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
# Row 76
$ws.Range("H76").Value = 29999
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
# Row 79
$ws.Range("H79").Value = 29999
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
# Row 126
$ws.Range("H126").Value = 2079.5557
$ws.Range("I126").Value = 2030.1818
$ws.Range("J126").Value = 2157.1428
$ws.Range("K126").Value = 6090.5454
$ws.Range("L126").Value = 6471.428400000001
$ws.Range("M126").Value = -3620.5454
$ws.Range("N126").Value = -11411.4284

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 19
$ws.Range("H19").Value = 2533.3333
$ws.Range("I19").Value = 2222.7273
$ws.Range("J19").Value = 5950
$ws.Range("K19").Value = 2222.7273
$ws.Range("L19").Value = 5950
$ws.Range("M19").Value = -2048.7273
$ws.Range("N19").Value = -6298
# Row 132
$ws.Range("H132").Value = 2336.848
$ws.Range("I132").Value = 1944.4333
$ws.Range("J132").Value = 3072.625
$ws.Range("K132").Value = 5833.2999
$ws.Range("L132").Value = 9217.875
$ws.Range("M132").Value = -3303.2999
$ws.Range("N132").Value = -14277.875

Write-Output "edits applied"
